# Insert a new row at position 24 (shifts existing rows 24-56 down to 25-57)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(24).Insert()

# Populate the new row 24 - same market/category metadata as the (now shifted)
# row below it, but with its own date / price figures.
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 44587
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 24000
$ws.Range("M24").Value = 23500
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Provincia de Diguillín"
$ws.Range("P24").Value = 940
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
